# Weekly data refresh: insert a new record as row 67 (Vega Central Mapocho
# de Santiago / Poroto granado) and push the existing rows 67-177 down to
# 68-178, matching the new A1:R178 used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 67 - this shifts every row
# at/after 67 down by one (old row 67 -> 68, ..., old row 177 -> 178).
$ws.Rows.Item(67).EntireRow.Insert()

# Populate the newly inserted row 67 with the new weekly observation.
$ws.Range("A67").Value = 9
$ws.Range("B67").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C67").Value = "Metropolitana"
$ws.Range("D67").Value = 44540
$ws.Range("E67").Value = 13
$ws.Range("F67").Value = 100112030
$ws.Range("G67").Value = "Poroto granado"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 43
$ws.Range("K67").Value = 36000
$ws.Range("L67").Value = 38000
$ws.Range("M67").Value = 37023
$ws.Range("N67").Value = "`$/saco 25 kilos"
$ws.Range("O67").Value = "Provincia de Limarí"
$ws.Range("P67").Value = 1481
$ws.Range("Q67").Value = 25
$ws.Range("R67").Value = "Hortaliza"
